# Apply the crypto price/volume refresh described by the commit diff.
# D column = Price, E column = Volume(1h) % change, both stored as plain text
# (not numbers) in the source workbook - values such as "615.30" or "1.00" must
# keep their literal formatting (trailing zeros) rather than being reinterpreted
# as numeric Doubles by Excel's usual text->number autoconversion.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.516.83'
$ws.Range("E2").Value = '  -1.72%  '
$ws.Range("D3").Value = '3.685.61'
$ws.Range("E3").Value = '  -2.88%  '
$ws.Range("E4").Value = '  +0.16%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '615.30'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  -0.24%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.02'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").Value = '3.679.54'
$ws.Range("E7").Value = '  -3.16%  '
$ws.Range("E8").Value = '  +0.13%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  -2.98%  '
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  -4.41%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.24'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  -2.52%  '
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.478'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  -5.00%  '
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.78'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  -2.45%  '
$ws.Range("E14").Value = '  -4.13%  '
$ws.Range("D15").Value = '4.301.28'
$ws.Range("E15").Value = '  -2.69%  '
$ws.Range("D16").Value = '3.687.26'
$ws.Range("E16").Value = '  -2.83%  '
$ws.Range("D17").Value = '69.515.54'
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("E18").Value = '  -2.00%  '
$ws.Range("E19").Value = '  -1.32%  '
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.26'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  -4.03%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '497.61'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -5.20%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.10'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  -3.90%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.719'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  -4.48%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.49'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -0.82%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.94'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  -2.66%  '
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.91'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  -5.78%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.12'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +1.29%  '
$ws.Range("E28").Value = '  +0.97%  '
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = $origStyle
$ws.Range("E30").Value = '  -3.31%  '
$ws.Range("E31").Value = '  -0.70%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.92'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  -0.78%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.03'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  -7.25%  '
$ws.Range("E34").Value = '  -2.20%  '
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  +0.11%  '
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.04'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  -1.70%  '
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.01'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  -3.28%  '
$ws.Range("E38").Value = '  +3.02%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.337'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -2.37%  '
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.91'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  -3.53%  '
$ws.Range("E41").Value = '  -7.30%  '
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.91'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  +2.85%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '425.41'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  -0.19%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.58'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  -1.88%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.54'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  -4.24%  '
$ws.Range("D46").Value = '2.930.73'
$ws.Range("E47").Value = '  -3.54%  '
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.26'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  -2.73%  '
$ws.Range("E49").Value = '  -0.07%  '
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.63'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  -3.60%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.44'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  -3.43%  '
